$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header C1: audioFalse -> currentPhase
$ws.Range("C1").Value = "currentPhase"

# Update C2 and C3 to share the new "train2P2" value (replacing the old
# per-row wav file names), consolidating two unique strings into one.
$ws.Range("C2").Value = "train2P2"
$ws.Range("C3").Value = "train2P2"
